# STAAS_Tagging.xlsx edit
#
# Commit intent: USER_NAME and API_TOKEN (API_KEY) are no longer read from
# this workbook -- they now come from environment variables in the script
# that consumes this sheet. So the "Fleet" sheet's USER_NAME / API_TOKEN
# columns (A, B) are removed, and the no-longer-used "Accessed_By" /
# Host / HostGroup / VM / Static column is removed as well, leaving just
# FUSION_SERVER, NAMESPACE and Application(/IT + App1..App32).

$wb = $excel.ActiveWorkbook

$fleet = $wb.Worksheets.Item("Fleet")

# Remove USER_NAME (A) and API_TOKEN (B) columns entirely.
$fleet.Columns("A:B").Delete()

# Remove the (now shifted-to-D) Accessed_By / Host / HostGroup / VM / Static
# column -- those values aren't used/tagged from this sheet any more.
$fleet.Columns("D:D").Delete()

# Keep the on-sheet cursor roughly where it was (shifted left one column
# since a column was removed ahead of the old selection), without stealing
# focus from the workbook's actual active tab.
$fleet.Range("E10").Select()

$tagging = $wb.Worksheets.Item("Tagging_map")
$tagging.Activate()
